$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 (n=10), "Enhanced DNC" block ---
$ws.Range("C14").Value = 0.0000481606
$ws.Range("D14").Value = 0.0000557899
$ws.Range("E14").Value = 0.0000286102
$ws.Range("F14").Value = 0.0000288486
$ws.Range("G14").Value = 0.0000238419
$ws.Range("H14").Value = 0.0000257492
$ws.Range("I14").Value = 0.0000195503
$ws.Range("J14").Value = 0.000023365
$ws.Range("K14").Value = 0.0000202656
$ws.Range("L14").Value = 0.0000193119
$ws.Range("M14").Value = 0.0000481606

# --- Row 15 (n=100) ---
$ws.Range("C15").Value = 0.0002086163
$ws.Range("D15").Value = 0.0003092289
$ws.Range("E15").Value = 0.0001916885
$ws.Range("F15").Value = 0.0001826286
$ws.Range("G15").Value = 0.0001907349
$ws.Range("H15").Value = 0.0001881123
$ws.Range("I15").Value = 0.0001854897
$ws.Range("J15").Value = 0.0001833439
$ws.Range("K15").Value = 0.0001802444
$ws.Range("L15").Value = 0.000187397
$ws.Range("M15").Value = 0.0002086163

# --- Row 16 (n=1000) ---
$ws.Range("C16").Value = 0.002692461
$ws.Range("D16").Value = 0.0021409988
$ws.Range("E16").Value = 0.0020980835
$ws.Range("F16").Value = 0.0028743744
$ws.Range("G16").Value = 0.002177
$ws.Range("H16").Value = 0.0021882057
$ws.Range("I16").Value = 0.0023241043
$ws.Range("J16").Value = 0.0022110939
$ws.Range("K16").Value = 0.0021512508
$ws.Range("L16").Value = 0.0021877289
$ws.Range("M16").Value = 0.002692461

# --- Row 17 (n=10000) ---
$ws.Range("C17").Value = 0.0294098854
$ws.Range("D17").Value = 0.0282216072
$ws.Range("E17").Value = 0.030138731
$ws.Range("F17").Value = 0.0300579071
$ws.Range("G17").Value = 0.0295789242
$ws.Range("H17").Value = 0.0299756527
$ws.Range("I17").Value = 0.0318930149
$ws.Range("J17").Value = 0.0303976536
$ws.Range("K17").Value = 0.0294623375
$ws.Range("L17").Value = 0.0288746357
$ws.Range("M17").Value = 0.0294098854

# --- Row 18 (n=100000) ---
$ws.Range("C18").Value = 0.5436584949
$ws.Range("D18").Value = 0.5398054123
$ws.Range("E18").Value = 0.6011793613
$ws.Range("F18").Value = 0.5864839554
$ws.Range("G18").Value = 0.6017947197
$ws.Range("H18").Value = 0.5902104378
$ws.Range("I18").Value = 0.5673451424
$ws.Range("J18").Value = 0.6004822254
$ws.Range("K18").Value = 0.573120594
$ws.Range("L18").Value = 0.6028115749
$ws.Range("M18").Value = 0.5436584949

# --- Update the selection to match the final saved cursor position ---
$ws.Range("M22").Select()
